# Auto-generated edit script applying the diff from the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("AB3").Value = 0.7058823529411765
$ws.Range("AC3").Value = 30700
$ws.Range("AD3").Value = 17
$ws.Range("AH3").Value = 0.29411764705882354
$ws.Range("AK3").Value = 0.7058823529411765
$ws.Range("AL3").Value = 31800
$ws.Range("AM3").Value = 17
$ws.Range("AQ3").Value = 0.29411764705882354
$ws.Range("AT3").Value = 0.7058823529411765
$ws.Range("AU3").Value = 33500
$ws.Range("AV3").Value = 16
$ws.Range("AZ3").Value = 0.3125
$ws.Range("B3").Value = 37900
$ws.Range("BC3").Value = 0.6875
$ws.Range("BD3").Value = 40300
$ws.Range("BE3").Value = 15
$ws.Range("BF3").Value = 5
$ws.Range("BG3").Value = 5
$ws.Range("BI3").Value = 0.3333333333333333
$ws.Range("BJ3").Value = 1
$ws.Range("BL3").Value = 0.6666666666666666
$ws.Range("BM3").Value = 34700
$ws.Range("BN3").Value = 16
$ws.Range("BR3").Value = 0.3125
$ws.Range("BU3").Value = 0.6875
$ws.Range("BV3").Value = 43200
$ws.Range("BW3").Value = 16
$ws.Range("C3").Value = 15
$ws.Range("CA3").Value = 0.3125
$ws.Range("CD3").Value = 0.6875
$ws.Range("CE3").Value = 31500
$ws.Range("CF3").Value = 17
$ws.Range("CJ3").Value = 0.29411764705882354
$ws.Range("CM3").Value = 0.7058823529411765
$ws.Range("CN3").Value = 35350
$ws.Range("CO3").Value = 16.2
$ws.Range("CP3").Value = 5
$ws.Range("CQ3").Value = 5
$ws.Range("CS3").Value = 0.3093137254901961
$ws.Range("CT3").Value = 1
$ws.Range("CV3").Value = 0.690686274509804
$ws.Range("CW3").Value = 30700
$ws.Range("CX3").Value = 17
$ws.Range("DB3").Value = 0.3333333333333333
$ws.Range("DE3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3333333333333333
$ws.Range("J3").Value = 0.6666666666666666
$ws.Range("K3").Value = 36000
$ws.Range("L3").Value = 16
$ws.Range("P3").Value = 0.3125
$ws.Range("S3").Value = 0.6875
$ws.Range("T3").Value = 33900
$ws.Range("U3").Value = 17
$ws.Range("Y3").Value = 0.29411764705882354
# Row 4
$ws.Range("AB4").Value = 0.9152542372881356
$ws.Range("AC4").Value = 7600
$ws.Range("AD4").Value = 52
$ws.Range("AH4").Value = 0.09615384615384616
$ws.Range("AK4").Value = 0.9038461538461539
$ws.Range("AL4").Value = 18800
$ws.Range("AM4").Value = 56
$ws.Range("AQ4").Value = 0.08928571428571429
$ws.Range("AT4").Value = 0.9107142857142857
$ws.Range("AU4").Value = 13000
$ws.Range("AV4").Value = 50
$ws.Range("AZ4").Value = 0.1
$ws.Range("B4").Value = 19700
$ws.Range("BC4").Value = 0.9
$ws.Range("BD4").Value = 13200
$ws.Range("BE4").Value = 59
$ws.Range("BF4").Value = 5
$ws.Range("BG4").Value = 5
$ws.Range("BI4").Value = 0.0847457627118644
$ws.Range("BJ4").Value = 1
$ws.Range("BL4").Value = 0.9152542372881356
$ws.Range("BM4").Value = 18600
$ws.Range("BN4").Value = 53
$ws.Range("BR4").Value = 0.09433962264150944
$ws.Range("BU4").Value = 0.9056603773584906
$ws.Range("BV4").Value = 20600
$ws.Range("BW4").Value = 54
$ws.Range("C4").Value = 56
$ws.Range("CA4").Value = 0.09259259259259259
$ws.Range("CD4").Value = 0.9074074074074074
$ws.Range("CE4").Value = 15200
$ws.Range("CF4").Value = 51
$ws.Range("CJ4").Value = 0.09803921568627451
$ws.Range("CM4").Value = 0.9019607843137255
$ws.Range("CN4").Value = 15750
$ws.Range("CO4").Value = 53.9
$ws.Range("CP4").Value = 5
$ws.Range("CQ4").Value = 5
$ws.Range("CS4").Value = 0.09312290473959106
$ws.Range("CT4").Value = 1
$ws.Range("CV4").Value = 0.906877095260409
$ws.Range("CW4").Value = 7600
$ws.Range("CX4").Value = 59
$ws.Range("G4").Value = 0.08928571428571429
$ws.Range("J4").Value = 0.9107142857142857
$ws.Range("K4").Value = 13300
$ws.Range("T4").Value = 17500
$ws.Range("U4").Value = 59
$ws.Range("Y4").Value = 0.0847457627118644
# Row 5
$ws.Range("AB5").Value = 0.5454545454545454
$ws.Range("AC5").Value = 18600
$ws.Range("AD5").Value = 16
$ws.Range("AE5").Value = 5
$ws.Range("AF5").Value = 5
$ws.Range("AH5").Value = 0.3125
$ws.Range("AI5").Value = 1
$ws.Range("AK5").Value = 0.6875
$ws.Range("AL5").Value = 16400
$ws.Range("AM5").Value = 15
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 4
$ws.Range("AQ5").Value = 0.26666666666666666
$ws.Range("AR5").Value = 0.8
$ws.Range("AT5").Value = 0.7333333333333333
$ws.Range("AU5").Value = 16100
$ws.Range("AV5").Value = 13
$ws.Range("AW5").Value = 5
$ws.Range("AX5").Value = 5
$ws.Range("AZ5").Value = 0.38461538461538464
$ws.Range("B5").Value = 19200
$ws.Range("BA5").Value = 1
$ws.Range("BC5").Value = 0.6153846153846154
$ws.Range("BD5").Value = 22000
$ws.Range("BE5").Value = 14
$ws.Range("BF5").Value = 4
$ws.Range("BG5").Value = 4
$ws.Range("BI5").Value = 0.2857142857142857
$ws.Range("BJ5").Value = 0.8
$ws.Range("BL5").Value = 0.7142857142857143
$ws.Range("BM5").Value = 14300
$ws.Range("BN5").Value = 10
$ws.Range("BO5").Value = 5
$ws.Range("BP5").Value = 5
$ws.Range("BR5").Value = 0.5
$ws.Range("BS5").Value = 1
$ws.Range("BU5").Value = 0.5
$ws.Range("BV5").Value = 12900
$ws.Range("BW5").Value = 10
$ws.Range("BX5").Value = 5
$ws.Range("BY5").Value = 5
$ws.Range("C5").Value = 9
$ws.Range("CA5").Value = 0.5
$ws.Range("CB5").Value = 1
$ws.Range("CD5").Value = 0.5
$ws.Range("CE5").Value = 14100
$ws.Range("CF5").Value = 14
$ws.Range("CG5").Value = 5
$ws.Range("CH5").Value = 5
$ws.Range("CJ5").Value = 0.35714285714285715
$ws.Range("CK5").Value = 1
$ws.Range("CM5").Value = 0.6428571428571429
$ws.Range("CN5").Value = 16480
$ws.Range("CO5").Value = 12.4
$ws.Range("CP5").Value = 4.7
$ws.Range("CQ5").Value = 4.7
$ws.Range("CS5").Value = 0.39500735375735374
$ws.Range("CT5").Value = 0.9399999999999998
$ws.Range("CV5").Value = 0.6049926462426463
$ws.Range("CW5").Value = 12900
$ws.Range("CX5").Value = 16
$ws.Range("CY5").Value = 5
$ws.Range("CZ5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("DB5").Value = 0.5555555555555556
$ws.Range("DC5").Value = 1
$ws.Range("DE5").Value = 0.4444444444444444
$ws.Range("E5").Value = 5
$ws.Range("G5").Value = 0.5555555555555556
$ws.Range("H5").Value = 1
$ws.Range("J5").Value = 0.4444444444444444
$ws.Range("K5").Value = 14300
$ws.Range("L5").Value = 12
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 4
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("Q5").Value = 0.8
$ws.Range("S5").Value = 0.6666666666666666
$ws.Range("T5").Value = 16900
$ws.Range("U5").Value = 11
$ws.Range("V5").Value = 5
$ws.Range("W5").Value = 5
$ws.Range("Y5").Value = 0.45454545454545453
$ws.Range("Z5").Value = 1
# Row 6
$ws.Range("AB6").Value = 0.9107142857142857
$ws.Range("AC6").Value = 6100
$ws.Range("AD6").Value = 64
$ws.Range("AH6").Value = 0.078125
$ws.Range("AK6").Value = 0.921875
$ws.Range("AL6").Value = 12900
$ws.Range("AM6").Value = 57
$ws.Range("AQ6").Value = 0.08771929824561403
$ws.Range("AT6").Value = 0.9122807017543859
$ws.Range("AU6").Value = 12900
$ws.Range("AV6").Value = 58
$ws.Range("AZ6").Value = 0.08620689655172414
$ws.Range("B6").Value = 34100
$ws.Range("BC6").Value = 0.9137931034482759
$ws.Range("BD6").Value = 16000
$ws.Range("BE6").Value = 59
$ws.Range("BF6").Value = 5
$ws.Range("BG6").Value = 5
$ws.Range("BI6").Value = 0.0847457627118644
$ws.Range("BJ6").Value = 1
$ws.Range("BL6").Value = 0.9152542372881356
$ws.Range("BM6").Value = 18300
$ws.Range("BN6").Value = 60
$ws.Range("BR6").Value = 0.08333333333333333
$ws.Range("BU6").Value = 0.9166666666666666
$ws.Range("BV6").Value = 9500
$ws.Range("BW6").Value = 60
$ws.Range("C6").Value = 55
$ws.Range("CA6").Value = 0.08333333333333333
$ws.Range("CD6").Value = 0.9166666666666666
$ws.Range("CE6").Value = 10900
$ws.Range("CF6").Value = 64
$ws.Range("CJ6").Value = 0.078125
$ws.Range("CM6").Value = 0.921875
$ws.Range("CN6").Value = 13890
$ws.Range("CO6").Value = 59.6
$ws.Range("CP6").Value = 5
$ws.Range("CQ6").Value = 5
$ws.Range("CS6").Value = 0.08411485087357538
$ws.Range("CT6").Value = 1
$ws.Range("CV6").Value = 0.9158851491264246
$ws.Range("CW6").Value = 6100
$ws.Range("CX6").Value = 64
$ws.Range("DB6").Value = 0.09090909090909091
$ws.Range("DE6").Value = 0.9090909090909091
$ws.Range("G6").Value = 0.09090909090909091
$ws.Range("J6").Value = 0.9090909090909091
$ws.Range("K6").Value = 9800
$ws.Range("L6").Value = 63
$ws.Range("P6").Value = 0.07936507936507936
$ws.Range("S6").Value = 0.9206349206349206
$ws.Range("T6").Value = 8400
$ws.Range("U6").Value = 56
$ws.Range("Y6").Value = 0.08928571428571429
